# Daily auto-update: prepend a new price row for 2026-02-05, pushing the
# existing rows (and the trailing 2025-11-21 row) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 2 (below the header row 1).
# This shifts all existing data rows down by one, growing the used range
# from A1:D77 to A1:D78.
$ws.Rows.Item(2).Insert()

# New row 2: same B/C/D figures as every other row, dated one day after
# the previous newest entry (2026-02-04 -> 2026-02-05).
# Force the date column to be stored as plain text (matching every other
# row in the sheet) instead of being auto-converted to a date serial.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-02-05"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
